$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.1835443037974684
$ws.Range("C2").Value = 0.5886075949367089
$ws.Range("J2").Value = 0.01265822784810127
$ws.Range("P2").Value = 0.1360759493670886
$ws.Range("S2").Value = 0.07911392405063292
$ws.Range("B3").Value = 0.007936507936507936
$ws.Range("C3").Value = 0.02116402116402116
$ws.Range("J3").Value = 0.01322751322751323
$ws.Range("O3").Value = 0.002645502645502645
$ws.Range("P3").Value = 0.7671957671957672
$ws.Range("S3").Value = 0.1878306878306878
$ws.Range("J4").Value = 0.02352941176470588
$ws.Range("P4").Value = 0.8117647058823529
$ws.Range("S4").Value = 0.1647058823529412
$ws.Range("P5").Value = 0.3333333333333333
$ws.Range("S5").Value = 0.6666666666666666
$ws.Range("B6").Value = 0.04631578947368421
$ws.Range("D6").Value = 0.01263157894736842
$ws.Range("F6").Value = 0.06736842105263158
$ws.Range("J6").Value = 0.271578947368421
$ws.Range("O6").Value = 0.02315789473684211
$ws.Range("Q6").Value = 0.1810526315789474
$ws.Range("R6").Value = 0.09473684210526316
$ws.Range("S6").Value = 0.3031578947368421
$ws.Range("B7").Value = 0.1206030150753769
$ws.Range("D7").Value = 0.02261306532663317
$ws.Range("F7").Value = 0.05276381909547739
$ws.Range("J7").Value = 0.1381909547738693
$ws.Range("O7").Value = 0.01758793969849246
$ws.Range("Q7").Value = 0.1884422110552764
$ws.Range("R7").Value = 0.0678391959798995
$ws.Range("S7").Value = 0.3919597989949749
$ws.Range("B8").Value = 0.1003236245954693
$ws.Range("D8").Value = 0.008629989212513484
$ws.Range("E8").Value = 0.002157497303128371
$ws.Range("F8").Value = 0.07011866235167206
$ws.Range("J8").Value = 0.1477885652642934
$ws.Range("O8").Value = 0.009708737864077669
$ws.Range("Q8").Value = 0.1909385113268608
$ws.Range("R8").Value = 0.1057173678532902
$ws.Range("S8").Value = 0.3646170442286947
$ws.Range("B9").Value = 0.1212871287128713
$ws.Range("D9").Value = 0.01237623762376238
$ws.Range("F9").Value = 0.0594059405940594
$ws.Range("J9").Value = 0.1262376237623762
$ws.Range("O9").Value = 0.01732673267326733
$ws.Range("Q9").Value = 0.1658415841584159
$ws.Range("R9").Value = 0.0891089108910891
$ws.Range("S9").Value = 0.4084158415841584
$ws.Range("B10").Value = 0.1097108969607116
$ws.Range("D10").Value = 0.02149740548554485
$ws.Range("E10").Value = 0.0007412898443291327
$ws.Range("F10").Value = 0.0737583395107487
$ws.Range("J10").Value = 0.1104521868050408
$ws.Range("O10").Value = 0.01334321719792439
$ws.Range("Q10").Value = 0.2138621200889548
$ws.Range("R10").Value = 0.09006671608598962
$ws.Range("S10").Value = 0.3665678280207561
$ws.Range("G11").Value = 0.1479674796747968
$ws.Range("J11").Value = 0.1024390243902439
$ws.Range("K11").Value = 0.1886178861788618
$ws.Range("L11").Value = 0.5544715447154471
$ws.Range("S11").Value = 0.006504065040650406
$ws.Range("G12").Value = 0.7321937321937322
$ws.Range("J12").Value = 0.207977207977208
$ws.Range("K12").Value = 0.008547008547008548
$ws.Range("L12").Value = 0.02849002849002849
$ws.Range("S12").Value = 0.02279202279202279
$ws.Range("G13").Value = 0.6739130434782609
$ws.Range("J13").Value = 0.2934782608695652
$ws.Range("S13").Value = 0.03260869565217391
$ws.Range("F15").Value = 0.01511879049676026
$ws.Range("H15").Value = 0.1749460043196544
$ws.Range("I15").Value = 0.06479481641468683
$ws.Range("J15").Value = 0.3542116630669546
$ws.Range("K15").Value = 0.06695464362850972
$ws.Range("M15").Value = 0.01295896328293736
$ws.Range("N15").Value = 0.002159827213822894
$ws.Range("O15").Value = 0.06263498920086392
$ws.Range("S15").Value = 0.2462203023758099
$ws.Range("F16").Value = 0.01627906976744186
$ws.Range("H16").Value = 0.1604651162790698
$ws.Range("I16").Value = 0.08372093023255814
$ws.Range("J16").Value = 0.3744186046511628
$ws.Range("K16").Value = 0.1279069767441861
$ws.Range("M16").Value = 0.02325581395348837
$ws.Range("O16").Value = 0.05348837209302326
$ws.Range("S16").Value = 0.1604651162790698
$ws.Range("F17").Value = 0.02361396303901437
$ws.Range("H17").Value = 0.1765913757700205
$ws.Range("I17").Value = 0.09240246406570841
$ws.Range("J17").Value = 0.4158110882956879
$ws.Range("K17").Value = 0.09650924024640657
$ws.Range("M17").Value = 0.02053388090349076
$ws.Range("O17").Value = 0.06160164271047228
$ws.Range("S17").Value = 0.1129363449691992
$ws.Range("F18").Value = 0.03555555555555556
$ws.Range("H18").Value = 0.1911111111111111
$ws.Range("I18").Value = 0.09777777777777778
$ws.Range("J18").Value = 0.4155555555555556
$ws.Range("K18").Value = 0.06888888888888889
$ws.Range("M18").Value = 0.01555555555555556
$ws.Range("O18").Value = 0.06666666666666667
$ws.Range("S18").Value = 0.1088888888888889
$ws.Range("F19").Value = 0.01396431342125679
$ws.Range("H19").Value = 0.2040341349883631
$ws.Range("I19").Value = 0.07874321179208689
$ws.Range("J19").Value = 0.3719937936384795
$ws.Range("K19").Value = 0.1078355314197052
$ws.Range("M19").Value = 0.02055857253685027
$ws.Range("N19").Value = 0.0003878975950349108
$ws.Range("O19").Value = 0.07408844065166796
$ws.Range("S19").Value = 0.1283941039565555
